# Edit script for ncsu_api_locations.xlsx
# Implements:
#  1. Sheet1: column K custom width removal (best-effort), sheetView selection change
#     (topLeftCell cleared, active cell -> A14)
#  2. Sheet1: split the old "notes" column L into:
#       - column L: a couple of rows now hold "actual start date: ..." notes, all
#         other previously-populated L cells become blank
#       - column M (new): holds the "#caution..." / "I can see on the more info..."
#         text that used to live in column L
#  3. New shared strings ("actual start date: 12/1/1870" / "actual start date: 1/1/1887")
#     get introduced automatically by simply writing that text into cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$CAUTION_TEXT = "#caution on website not to use for pcp"
$CHRONOS_TEXT = "I can see on the more info tab they have todays data, but no data available on api or through chronos"

# Rows whose column-L "caution" note moves to column M, with the s=3 (light) style.
$rowsCautionS3 = @(29,30,31,32,33,34,35,36,37,38,39,40,41,42,43,44,46,69,70,79,80,81,82,83,84,85,86,87)

# Rows whose column-L "chronos" note moves to column M, with the s=5 (dark) style.
$rowsChronosS5 = @(45,47,48,50,51,52,53,54,55,56,65,67,68)

foreach ($r in $rowsCautionS3) {
    $lcell = $ws.Range("L$r")
    $mcell = $ws.Range("M$r")
    $lcell.Copy() | Out-Null
    $mcell.PasteSpecial(-4122) | Out-Null   # xlPasteFormats
    $mcell.Value2 = $CAUTION_TEXT
    $lcell.ClearContents() | Out-Null
}

foreach ($r in $rowsChronosS5) {
    $lcell = $ws.Range("L$r")
    $mcell = $ws.Range("M$r")
    $lcell.Copy() | Out-Null
    $mcell.PasteSpecial(-4122) | Out-Null   # xlPasteFormats
    $mcell.Value2 = $CHRONOS_TEXT
    $lcell.ClearContents() | Out-Null
}

# Special rows: column L keeps a value but it changes from the "chronos" note to a
# brand-new "actual start date" note; the old chronos note moves to the new column M.
$lcell = $ws.Range("L49")
$mcell = $ws.Range("M49")
$lcell.Copy() | Out-Null
$mcell.PasteSpecial(-4122) | Out-Null
$mcell.Value2 = $CHRONOS_TEXT
$lcell.Value2 = "actual start date: 12/1/1870"

$lcell = $ws.Range("L66")
$mcell = $ws.Range("M66")
$lcell.Copy() | Out-Null
$mcell.PasteSpecial(-4122) | Out-Null
$mcell.Value2 = $CHRONOS_TEXT
$lcell.Value2 = "actual start date: 1/1/1887"

# Cosmetic worksheet-level tweaks -------------------------------------------------

# Update the saved selection/view: no frozen top-left cell, active cell now A14.
$ws.Range("A14").Select() | Out-Null

Write-Output "done"
